# config.xlsx edit:
#  - add the DMS home-site URL into the "DMSHomeUrl" row (B19) of the
#    ROBOT PARAMETERS table, pulling in a new shared string
#  - move the active selection to B20 (reflecting where the user ended
#    up after filling in the new value)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 = DMSHomeUrl ; give it its Value (column B) -> the DMS home url
$ws.Cells.Item(19, 2).Value = "https://defradev.sharepoint.com/sites/EADMSRoboticsHomeSite"

# Update the active cell / selection to B20
$ws.Range("B20").Select()
